# Daily attendance processing - 2026-01-07 07:45:23
# Reorders the "Recorded By" list in column G so that "System" (exact case)
# is moved to the front of the comma-separated list, preserving the
# relative order of the remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value()

    if ($null -eq $value -or $value -eq "") {
        continue
    }

    $parts = $value -split ",\s*"
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $hasSystem = $false
    foreach ($p in $trimmed) {
        if ($p.Equals("System")) {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $rest = @()
        foreach ($p in $trimmed) {
            if (-not $p.Equals("System")) {
                $rest += $p
            }
        }
        $newParts = @("System") + $rest
        $newValue = [string]::Join(", ", $newParts)

        if (-not $newValue.Equals($value)) {
            $cell.Value = $newValue
        }
    }
}
